# issue #5: add legislator_id, name, date into dataframe
#
# The stock ("股票") sheet gains three new trailing columns: date,
# legislator_name, legislator_id. Populate the header row with the new
# column names (matching the bold / centered / bordered look already used
# by the existing header cells B1:G1) and fill every data row with the
# report date, the legislator's name and her numeric id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "田秋堇"
$legislatorId = 1316
$reportDate = "2012-04-10"

# --- header row (row 1) --------------------------------------------------
$headers = @{ "H1" = "date"; "I1" = "legislator_name"; "J1" = "legislator_id" }
foreach ($addr in $headers.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $headers[$addr]
    # Reproduce the look of the existing header cells (bold, centered,
    # top-aligned, thin border all round) instead of leaving General style.
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# --- data rows (rows 2-4) -------------------------------------------------
for ($r = 2; $r -le 4; $r++) {
    $dateCell = $ws.Range("H" + $r)
    # Force text storage so "2012-04-10" is kept as the literal string it
    # is, instead of being auto-converted into an Excel date serial number.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $reportDate

    $ws.Range("I" + $r).Value = $legislatorName
    $ws.Range("J" + $r).Value = $legislatorId
}
